# Regenerate the "K" column (column G) values in the save_data sheet.
# This mirrors a re-run of the data pipeline that recalculates the K stat
# (formerly "Strike#") for every recorded game row (rows 2-48), after which
# std/mean and s_vals were recomputed upstream. Only the stored K values
# change in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new K value pairs taken from the regenerated save data.
$kUpdates = @(
    @(2, 1), @(3, 0), @(4, 1), @(5, 1), @(6, 0), @(7, 1), @(8, 3), @(9, 0),
    @(10, 3), @(11, 1), @(12, 0), @(13, 2), @(14, 2), @(15, 0), @(16, 1),
    @(17, 1), @(18, 3), @(19, 0), @(20, 1), @(21, 1), @(22, 0), @(23, 0),
    @(24, 1), @(25, 2), @(26, 0), @(27, 0), @(28, 1), @(29, 0), @(30, 0),
    @(31, 0), @(32, 3), @(33, 0), @(34, 1), @(35, 0), @(36, 1), @(37, 1),
    @(38, 2), @(39, 2), @(40, 0), @(41, 1), @(42, 0), @(43, 2), @(44, 4),
    @(45, 2), @(46, 2), @(47, 1), @(48, 1)
)

foreach ($update in $kUpdates) {
    $row = $update[0]
    $value = $update[1]
    $ws.Cells.Item($row, 7).Value = $value
}
